$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.7
$ws.Range("I3").Value = 2.5
$ws.Range("J3").Value = 3.15
$ws.Range("L3").Value = 3.05
$ws.Range("M3").Value = 1.02
$ws.Range("N3").Value = 9.7
$ws.Range("W3").Value = 9.5
$ws.Range("X3").Value = 14.5
$ws.Range("Y3").Value = 9.75
$ws.Range("Z3").Value = 32
$ws.Range("AA3").Value = 22
$ws.Range("AB3").Value = 28
$ws.Range("AG3").Value = 8.25
$ws.Range("AH3").Value = 12.5
$ws.Range("AI3").Value = 9.5
$ws.Range("AJ3").Value = 28
$ws.Range("AK3").Value = 21
$ws.Range("AL3").Value = 29
$ws.Range("AN3").Value = 4.7
$ws.Range("AO3").Value = 14
$ws.Range("AW3").Value = 4.5
$ws.Range("AX3").Value = 13
$ws.Range("AY3").Value = 19.5
$ws.Range("AZ3").Value = 55
$ws.Range("BA3").Value = 80
$ws.Range("BB3").Value = 200

# Row 4
$ws.Range("G4").Value = 1.98
$ws.Range("I4").Value = 4.05
$ws.Range("P4").Value = 2.72
$ws.Range("V4").Value = 1.82
$ws.Range("W4").Value = 6.3
$ws.Range("X4").Value = 9
$ws.Range("Z4").Value = 18
$ws.Range("AB4").Value = 30
$ws.Range("AC4").Value = 8
$ws.Range("AD4").Value = 5.9
$ws.Range("AG4").Value = 10.75
$ws.Range("AH4").Value = 23
$ws.Range("AJ4").Value = 70
$ws.Range("AN4").Value = 3.8
$ws.Range("AP4").Value = 18
$ws.Range("AR4").Value = 70
$ws.Range("AY4").Value = 25

# Row 6
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11

# Row 7
$ws.Range("Q7").Value = 1.9
$ws.Range("R7").Value = 1.9

# Row 8
$ws.Range("G8").Value = 3.6
$ws.Range("I8").Value = 2.05
$ws.Range("N8").Value = 8
$ws.Range("AF8").Value = 67
$ws.Range("AH8").Value = 9
$ws.Range("AQ8").Value = 81
